# Apply the edit: insert one new data row before row 153, shifting the
# existing rows 153:228 down to 154:229, and populate the newly inserted
# row 153 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 153; this pushes old row153..row228 down
# to row154..row229 and keeps all their existing content/formatting intact.
$ws.Rows(153).Insert()

# Fill in the values for the newly inserted row 153.
$ws.Cells.Item(153, 1).Value = 7
$ws.Cells.Item(153, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(153, 3).Value = "Ñuble"
$ws.Cells.Item(153, 4).Value = 44813
$ws.Cells.Item(153, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(153, 5).Value = 16
$ws.Cells.Item(153, 6).Value = "Fruta"
$ws.Cells.Item(153, 7).Value = 100104
$ws.Cells.Item(153, 8).Value = "Frutos de pepita"
$ws.Cells.Item(153, 9).Value = 100104005
$ws.Cells.Item(153, 10).Value = "Pera"
$ws.Cells.Item(153, 11).Value = "Packham's Triumph"
$ws.Cells.Item(153, 12).Value = "Primera"
$ws.Cells.Item(153, 13).Value = 120
$ws.Cells.Item(153, 14).Value = 11000
$ws.Cells.Item(153, 15).Value = 12000
$ws.Cells.Item(153, 16).Value = 11500
$ws.Cells.Item(153, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(153, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(153, 19).Value = 719
$ws.Cells.Item(153, 20).Value = 16
